$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Comparison")

# Table 1 (Matrix size / Binary / Ladder / Exponential) rows 3-15
$ws.Range("B3").Value = 105.0
$ws.Range("C3").Value = 1996.0
$ws.Range("D3").Value = 3495.0
$ws.Range("B4").Value = 119.0
$ws.Range("C4").Value = 1912.0
$ws.Range("D4").Value = 4299.0
$ws.Range("B5").Value = 166.0
$ws.Range("C5").Value = 2099.0
$ws.Range("D5").Value = 4821.0
$ws.Range("B6").Value = 312.0
$ws.Range("C6").Value = 2356.0
$ws.Range("D6").Value = 5524.0
$ws.Range("B7").Value = 602.0
$ws.Range("C7").Value = 2595.0
$ws.Range("D7").Value = 6344.0
$ws.Range("B8").Value = 1214.0
$ws.Range("C8").Value = 2788.0
$ws.Range("D8").Value = 6608.0
$ws.Range("B9").Value = 2198.0
$ws.Range("C9").Value = 2603.0
$ws.Range("D9").Value = 6491.0
$ws.Range("B10").Value = 4916.0
$ws.Range("C10").Value = 3122.0
$ws.Range("D10").Value = 7141.0
$ws.Range("B11").Value = 13962.0
$ws.Range("C11").Value = 5222.0
$ws.Range("D11").Value = 8653.0
$ws.Range("B12").Value = 29901.0
$ws.Range("C12").Value = 8049.0
$ws.Range("D12").Value = 11597.0
$ws.Range("B13").Value = 125517.0
$ws.Range("C13").Value = 14339.0
$ws.Range("D13").Value = 19382.0
$ws.Range("B14").Value = 252189.0
$ws.Range("C14").Value = 27870.0
$ws.Range("D14").Value = 39279.0
$ws.Range("B15").Value = 536612.0
$ws.Range("C15").Value = 46732.0
$ws.Range("D15").Value = 62442.0

# Table 2 rows 19-31
$ws.Range("B19").Value = 80.0
$ws.Range("C19").Value = 2421.0
$ws.Range("D19").Value = 5601.0
$ws.Range("B20").Value = 97.0
$ws.Range("C20").Value = 2358.0
$ws.Range("D20").Value = 5464.0
$ws.Range("B21").Value = 163.0
$ws.Range("C21").Value = 2386.0
$ws.Range("D21").Value = 5544.0
$ws.Range("B22").Value = 313.0
$ws.Range("C22").Value = 2421.0
$ws.Range("D22").Value = 6786.0
$ws.Range("B23").Value = 583.0
$ws.Range("C23").Value = 2449.0
$ws.Range("D23").Value = 5848.0
$ws.Range("B24").Value = 1119.0
$ws.Range("C24").Value = 2490.0
$ws.Range("D24").Value = 6328.0
$ws.Range("B25").Value = 2265.0
$ws.Range("C25").Value = 2692.0
$ws.Range("D25").Value = 6506.0
$ws.Range("B26").Value = 5240.0
$ws.Range("C26").Value = 3201.0
$ws.Range("D26").Value = 6990.0
$ws.Range("B27").Value = 13358.0
$ws.Range("C27").Value = 4320.0
$ws.Range("D27").Value = 7192.0
$ws.Range("B28").Value = 30603.0
$ws.Range("C28").Value = 6895.0
$ws.Range("D28").Value = 8200.0
$ws.Range("B29").Value = 62899.0
$ws.Range("C29").Value = 11338.0
$ws.Range("D29").Value = 9268.0
$ws.Range("B30").Value = 124890.0
$ws.Range("C30").Value = 19528.0
$ws.Range("D30").Value = 12612.0
$ws.Range("B31").Value = 256699.0
$ws.Range("C31").Value = 35854.0
$ws.Range("D31").Value = 15407.0

# Table 3 (ratio) rows 35-47
$ws.Range("B35").Value = 1.3125
$ws.Range("C35").Value = 0.8244527054935977
$ws.Range("D35").Value = 0.6239957150508838
$ws.Range("B36").Value = 1.2268041237113403
$ws.Range("C36").Value = 0.8108566581849025
$ws.Range("D36").Value = 0.7867862371888726
$ws.Range("B37").Value = 1.01840490797546
$ws.Range("C37").Value = 0.8797150041911148
$ws.Range("D37").Value = 0.8695887445887446
$ws.Range("B38").Value = 0.9968051118210862
$ws.Range("C38").Value = 0.973151590251962
$ws.Range("D38").Value = 0.8140288829944002
$ws.Range("B39").Value = 1.032590051457976
$ws.Range("C39").Value = 1.059616169865251
$ws.Range("D39").Value = 1.084815321477428
$ws.Range("B40").Value = 1.0848972296693475
$ws.Range("C40").Value = 1.1196787148594378
$ws.Range("D40").Value = 1.0442477876106195
$ws.Range("B41").Value = 0.9704194260485651
$ws.Range("C41").Value = 0.9669390787518574
$ws.Range("D41").Value = 0.9976944359053181
$ws.Range("B42").Value = 0.9381679389312977
$ws.Range("C42").Value = 0.9753202124336144
$ws.Range("D42").Value = 1.0216022889842633
$ws.Range("B43").Value = 1.045216349752957
$ws.Range("C43").Value = 1.2087962962962964
$ws.Range("D43").Value = 1.203142380422692
$ws.Range("B44").Value = 0.9770610724438781
$ws.Range("C44").Value = 1.1673676577229877
$ws.Range("D44").Value = 1.4142682926829269
$ws.Range("B45").Value = 1.9955325203898313
$ws.Range("C45").Value = 1.2646851296524961
$ws.Range("D45").Value = 2.091281829952525
$ws.Range("B46").Value = 2.019288974297382
$ws.Range("C46").Value = 1.427181482998771
$ws.Range("D46").Value = 3.1144148430066605
$ws.Range("B47").Value = 2.090432763664837
$ws.Range("C47").Value = 1.3033971105037094
$ws.Range("D47").Value = 4.052833127799053

$wb.Save()
